# Replace the Maltaspor fantasy-roster table (player / position / team)
# with an updated 2024 player list, keeping the existing header row and
# its formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - Oyuncu Adı (Player Name)
$colA = @(
    "Oyuncu Adı",
    "Cade Cunningham",
    "Derrick White",
    "Quentin Grimes",
    "Brandon Ingram",
    "Brandon Miller",
    "Herbert Jones",
    "Isaiah Hartenstein",
    "Bam Adebayo",
    "Kentavious Caldwell-Pope",
    "Damian Lillard",
    "Harrison Barnes",
    "Julius Randle",
    "Malik Monk",
    "Anthony Davis",
    "Cameron Johnson",
    "LaMelo Ball",
    "Nick Richards"
)

# Column B - Pozisyon (Position)
$colB = @(
    "Pozisyon",
    "PG,SG",
    "SG,SF",
    "SG,SF,PF",
    "SF,PF",
    "SG,SF",
    "SF,PF",
    "C",
    "C",
    "SG,SF",
    "PG",
    "SF,PF",
    "PF",
    "SG,SF",
    "PF,C",
    "SF,PF",
    "PG,SG",
    "C"
)

# Column C - Takım (Team)
$colC = @(
    "Takım",
    "Detroit Pistons",
    "Boston Celtics",
    "Dallas Mavericks",
    "New Orleans Pelicans",
    "Charlotte Hornets",
    "New Orleans Pelicans",
    "Oklahoma City Thunder",
    "Miami Heat",
    "Orlando Magic",
    "Milwaukee Bucks",
    "San Antonio Spurs",
    "Minnesota Timberwolves",
    "Sacramento Kings",
    "Los Angeles Lakers",
    "Brooklyn Nets",
    "Charlotte Hornets",
    "Charlotte Hornets"
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}
